# Applies the "1401ME79" marksheet edit: updates the score summary (rows
# 10-12), recomputes which questions were answered correctly (col A vs col
# B in rows 16-40), and drops the extra "attempt" pairs (cols D/E except
# rows 17-18, and cols G/H everywhere) that no longer apply.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Score summary block (rows 10-12)
# ---------------------------------------------------------------------

# A10/A11/A12 gain the bold "mtitleStyle" look already used by A9 - copy
# its formatting across instead of assigning a named style so the
# existing cellXfs entry (s="4") is reused rather than a new one minted.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 72
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = "72/112"

# ---------------------------------------------------------------------
# Per-question answer grid (rows 16-40)
# ---------------------------------------------------------------------

# Rows whose "Student Ans" (col A) now matches the "Correct Ans" (col B),
# i.e. the student answered correctly - give col A the green
# "correctStyle" look (copied from B10, which already carries s="5") and
# fill in the matching option text.
$correctRows = @(16,17,18,19,22,23,25,27,28,29,30,32,33,35,38,39)

$ws.Range("B10").Copy()
foreach ($r in $correctRows) {
    $ws.Range("A$r").PasteSpecial(-4122)
}

$ws.Range("A16").Value = "Option A"
$ws.Range("A17").Value = "Option D"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("A25").Value = "Option A"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A35").Value = "Option D"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"

# Rows 17 & 18 also keep a second question pair in D/E; D gets the same
# "answered correctly" treatment as col A above.
$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

# Drop the now-unused 3rd question-pair columns (G/H) entirely, and the
# 2nd question-pair columns (D/E) everywhere except rows 17-18.
$ws.Range("G15:H40").Clear()
$ws.Range("D16:E16").Clear()
$ws.Range("D19:E40").Clear()
